$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Add the new "metadata" sheet right after "data"
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Match page margins used by the "data" sheet (0.75/0.75/1/1/0.5/0.5 in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row (B1:G1), styled like the "data" sheet header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").PasteSpecial(-4122)

# Data row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Ultra-rare undescribed monogenic disorders"
$ws.Range("C2").Value = 195

# D2 must stay textual ("1.2"), not be coerced into the number 1.2
$ws.Range("D2").Formula = "=""1.2"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("E2").Value = "2020-05-07T14:27:15.216588Z"
$ws.Range("F2").Value = "2021-10-05 14:23:00.551278"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/195/?format=json"

# Restore "data" as the active sheet/selection (matches the unchanged bookViews)
$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null
